$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New translation rows to append (language, label, translation)
$rows = @(
    @("cs", "lab.liquid.tooltip.create", "Nový liquid"),
    @("cs", "lab.liquid.tooltip.create", "Vytvořit liquid"),
    @("cs", "lab.liquid.create.title", "Nový liquid"),
    @("cs", "lab.liquid.create.subtitle", "Liquidy je možné použít pro míchání mixů k vapování."),
    @("cs", "lab.liquid.name.label", "Jméno"),
    @("cs", "lab.liquid.description.label", "Popis"),
    @("cs", "lab.liquid.vendorId.label", "Výrobce"),
    @("cs", "lab.liquid.pg.label", "PG"),
    @("cs", "lab.liquid.vg.label", "VG"),
    @("cs", "lab.liquid.create.submit", "Vytvořit liquid"),
    @("cs", "lab.liquid.volume.label", "Objem"),
    @("cs", "lab.liquid.create.success", "Liquid [{{data.name}}] byl uložen.")
)

$startRow = 333
$templateRow = $startRow - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]

    # Reuse the formatting of the previous (last existing) row instead of
    # applying a named style, so the same cellXf is shared rather than a
    # new (duplicate) one being created.
    $ws.Range($ws.Cells.Item($templateRow, 1), $ws.Cells.Item($templateRow, 3)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3)).PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = $false

# Update the visible window / selection to reflect scrolling to the new rows
$ws.Application.ActiveWindow.ScrollRow = 328
$ws.Range("B342").Select()
